$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.224535333333333
$ws.Range("H2").Value = 12.673606
$ws.Range("I2").Value = 0.7043225486309714
$ws.Range("J2").Value = 0.7043225486309715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.345835
$ws.Range("N2").Value = 91.03750500000001
$ws.Range("O2").Value = 0.8527782452855476
$ws.Range("P2").Value = 0.8527782452855475
$ws.Range("Q2").Value = 128.1970521770033
$ws.Range("R2").Value = 1153.77346959303
$ws.Range("S2").Value = 0.6006309471365646
$ws.Range("T2").Value = 0.6006309471365646

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.224535333333333
$ws.Range("H3").Value = 12.673606
$ws.Range("I3").Value = 0.7043225486309714
$ws.Range("J3").Value = 0.7043225486309715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.224107666666666
$ws.Range("N3").Value = 9.672322999999999
$ws.Range("O3").Value = 0.09060383010029813
$ws.Range("P3").Value = 0.09060383010029811
$ws.Range("Q3").Value = 13.62035675630422
$ws.Range("R3").Value = 122.583210806738
$ws.Range("S3").Value = 0.0638143205319695
$ws.Range("T3").Value = 0.0638143205319695

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.224535333333333
$ws.Range("H4").Value = 12.673606
$ws.Range("I4").Value = 0.7043225486309714
$ws.Range("J4").Value = 0.7043225486309715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.014730333333333
$ws.Range("N4").Value = 6.044191
$ws.Range("O4").Value = 0.05661792461415433
$ws.Range("P4").Value = 0.05661792461415433
$ws.Range("Q4").Value = 8.511299480305109
$ws.Range("R4").Value = 76.601695322746
$ws.Range("S4").Value = 0.03987728096243739
$ws.Range("T4").Value = 0.03987728096243739

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.773477
$ws.Range("H5").Value = 5.320431
$ws.Range("I5").Value = 0.2956774513690286
$ws.Range("J5").Value = 0.2956774513690286
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.345835
$ws.Range("N5").Value = 91.03750500000001
$ws.Range("O5").Value = 0.8527782452855476
$ws.Range("P5").Value = 0.8527782452855475
$ws.Range("Q5").Value = 53.81764041829501
$ws.Range("R5").Value = 484.3587637646551
$ws.Range("S5").Value = 0.252147298148983
$ws.Range("T5").Value = 0.252147298148983

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.773477
$ws.Range("H6").Value = 5.320431
$ws.Range("I6").Value = 0.2956774513690286
$ws.Range("J6").Value = 0.2956774513690286
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.224107666666666
$ws.Range("N6").Value = 9.672322999999999
$ws.Range("O6").Value = 0.09060383010029813
$ws.Range("P6").Value = 0.09060383010029811
$ws.Range("Q6").Value = 5.717880792357
$ws.Range("R6").Value = 51.460927131213
$ws.Range("S6").Value = 0.02678950956832863
$ws.Range("T6").Value = 0.02678950956832862

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.773477
$ws.Range("H7").Value = 5.320431
$ws.Range("I7").Value = 0.2956774513690286
$ws.Range("J7").Value = 0.2956774513690286
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.014730333333333
$ws.Range("N7").Value = 6.044191
$ws.Range("O7").Value = 0.05661792461415433
$ws.Range("P7").Value = 0.05661792461415433
$ws.Range("Q7").Value = 3.573077907369
$ws.Range("R7").Value = 32.157701166321
$ws.Range("S7").Value = 0.01674064365171694
$ws.Range("T7").Value = 0.01674064365171694
